$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10192
$ws.Range("I32").Value = 9249.666999999999
$ws.Range("K32").Value = 9249.666999999999
$ws.Range("M32").Value = -8923.666999999999
$ws.Range("H33").Value = 1135
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5458
$ws.Range("H40").Value = 7758
$ws.Range("J40").Value = 8346.434999999999
$ws.Range("L40").Value = 8346.434999999999
$ws.Range("N40").Value = -8696.434999999999
$ws.Range("H58").Value = 4503
$ws.Range("I58").Value = 15
$ws.Range("J58").Value = 5625
$ws.Range("K58").Value = 45
$ws.Range("L58").Value = 16875
$ws.Range("M58").Value = 105
$ws.Range("N58").Value = -17175
$ws.Range("H64").Value = 20004000
$ws.Range("I64").Value = 50002500
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 50002500
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -50002252
$ws.Range("N64").Value = -5496
$ws.Range("H67").Value = 20004000
$ws.Range("I67").Value = 50002500
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 50002500
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -50001642
$ws.Range("N67").Value = -6716
$ws.Range("H74").Value = 8066.6665
$ws.Range("I74").Value = 7200
$ws.Range("J74").Value = 8500
$ws.Range("K74").Value = 7200
$ws.Range("L74").Value = 8500
$ws.Range("M74").Value = -6264
$ws.Range("N74").Value = -10372
$ws.Range("H77").Value = 8066.6665
$ws.Range("I77").Value = 7200
$ws.Range("J77").Value = 8500
$ws.Range("K77").Value = 36000
$ws.Range("L77").Value = 42500
$ws.Range("M77").Value = -31320
$ws.Range("N77").Value = -51860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 2500
$ws.Range("K45").Value = 2500
$ws.Range("M45").Value = -2123
$ws.Range("H97").Value = 574
$ws.Range("I97").Value = 574
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 574
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -78
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 4964.6665
$ws.Range("J122").Value = 4997.5
$ws.Range("L122").Value = 14992.5
$ws.Range("N122").Value = -19892.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 1000
$ws.Range("J23").Value = 1000
$ws.Range("L23").Value = 1000
$ws.Range("N23").Value = -1566
$ws.Range("H94").Value = 1177.1666
$ws.Range("I94").Value = 1177.1666
$ws.Range("K94").Value = 1177.1666
$ws.Range("M94").Value = -726.1666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 78947570
$ws.Range("I7").Value = 38461776
$ws.Range("J7").Value = 166666770
$ws.Range("K7").Value = 38461776
$ws.Range("L7").Value = 166666770
$ws.Range("M7").Value = -38461663
$ws.Range("N7").Value = -166666996
$ws.Range("H16").Value = 2833.3333
$ws.Range("I16").Value = 1750
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1750
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -1463
$ws.Range("N16").Value = -5574
$ws.Range("H62").Value = 41674916
$ws.Range("I62").Value = 50007896
$ws.Range("K62").Value = 50007896
$ws.Range("M62").Value = -50007272
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 41674916
$ws.Range("I65").Value = 50007896
$ws.Range("K65").Value = 250039480
$ws.Range("M65").Value = -250036360
$ws.Range("N65").ClearContents()
$ws.Range("H99").Value = 953.6
$ws.Range("I99").Value = 953.6
$ws.Range("K99").Value = 953.6
$ws.Range("M99").Value = 544.4
$ws.Range("H113").Value = 2833.3333
$ws.Range("I113").Value = 1750
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1750
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 420
$ws.Range("N113").Value = -9340
$ws.Range("H124").Value = 50326
$ws.Range("J124").Value = 50326
$ws.Range("L124").Value = 50326
$ws.Range("N124").Value = -55236
$ws.Range("H126").Value = 953.6
$ws.Range("I126").Value = 953.6
$ws.Range("K126").Value = 2860.8
$ws.Range("M126").Value = -390.8000000000002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1633.3334
$ws.Range("I9").Value = 1633.3334
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1633.3334
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1463.3334
$ws.Range("N9").ClearContents()
$ws.Range("H18").Value = 2500
$ws.Range("J18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("N18").Value = -3086
$ws.Range("H80").Value = 3043.125
$ws.Range("I80").Value = 3056.4285
$ws.Range("J80").Value = 2950
$ws.Range("K80").Value = 3056.4285
$ws.Range("L80").Value = 2950
$ws.Range("M80").Value = -2058.4285
$ws.Range("N80").Value = -4946
$ws.Range("H83").Value = 3043.125
$ws.Range("I83").Value = 3056.4285
$ws.Range("J83").Value = 2950
$ws.Range("K83").Value = 15282.1425
$ws.Range("L83").Value = 14750
$ws.Range("M83").Value = -10290.1425
$ws.Range("N83").Value = -24734

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 40000
$ws.Range("I93").Value = 40000
$ws.Range("K93").Value = 40000
$ws.Range("M93").Value = -38752
$ws.Range("H125").Value = 70664
$ws.Range("J125").Value = 70664
$ws.Range("L125").Value = 70664
$ws.Range("N125").Value = -80504
$ws.Range("H136").Value = 5052
$ws.Range("I136").Value = 4785.778
$ws.Range("K136").Value = 14357.334
$ws.Range("M136").Value = -11807.334
$ws.Range("N136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4256.857
$ws.Range("I96").Value = 3633
$ws.Range("J96").Value = 8000
$ws.Range("K96").Value = 3633
$ws.Range("L96").Value = 8000
$ws.Range("M96").Value = -2260
$ws.Range("N96").Value = -10746
$ws.Range("H136").Value = 4693.4443
$ws.Range("I136").Value = 4050.6
$ws.Range("J136").Value = 5497
$ws.Range("K136").Value = 12151.8
$ws.Range("L136").Value = 16491
$ws.Range("M136").Value = -9601.799999999999
$ws.Range("N136").Value = -21591
